# Generate Report for Archive
#
# Updates the localization status for two files (1d7a6d2c-...md and
# 2a92bb2d-...md) from "Ready for handoff" to "In Translation":
#   - Overview sheet: zh-cn (E) and de-de (F) status columns, rows 3 & 4
#   - zh-cn sheet: Status column (C), rows 3 & 4
#   - de-de sheet: Status column (C), rows 3 & 4
#
# The third file (2d529088-...md, row 5) keeps its "Ready for handoff"
# status and is left untouched.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
